$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the text of B16 (3.5.1.1 indicator passport, "Методы сбора данных" -> second row) ---
# The reporting form number/title changed from No.11 "Отчет о наркологических расстройствах"
# to No.10 "Отчет о психических и наркологических расстройствах".
$newText = "Минздрав КР, министерства и ведомства представляют сводный отчет в разрезе областей и районов по форме государственной статистической отчетности №10 «Отчет о психических и наркологических расстройствах»."
$target = $ws.Range("B16")
$target.Value = $newText

# --- Give the edited cell its own (new) font/style, distinguishing it from the other
#     data-entry cells in column B which keep sharing the original font/style. ---
$target.Font.Name = "Calibri"
$target.Font.Size = 11
$target.WrapText = $true
$target.VerticalAlignment = -4160
$target.Locked = $false

# --- Reflect where the user ended up working: scrolled down to row 13, with B16 selected. ---
$ws.Activate()
$target.Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
